$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A18").Value = "GFG"
$ws.Range("B18").Value = "Merge Sort For Linked lists.[Very Important]"
$ws.Range("E17").Select()
